$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 data
$ws.Range("A12").Value = 45708
$ws.Range("A12").NumberFormat = "yyyy-mm-dd"
$ws.Range("C12").Value = "Corte Adulto"
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = "julian"
$ws.Range("G12").Value = "Efectivo"
